# "UEN ACRA API fix"
#
# The VAT ID column (H) previously stored a hard-coded placeholder value
# (12345678) for every company. It is fixed to instead mirror the UEN
# column (A) via a formula, and two additional companies (with their ACRA
# UEN-derived data) are appended to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix column H (VAT ID) to reference the UEN in column A ---------------
# H2 becomes its own formula; H3:H5 becomes a shared formula group (set as
# one range assignment so Excel groups them together).
$ws.Range("H2").Formula = "=A2"
$ws.Range("H3:H5").Formula = "=A3"

# --- Apply thousands-separator number formatting to the financial columns -
$ws.Range("C2:E5").NumberFormat = "#,##0"

# Spacer row between the existing data and the newly appended rows.
$ws.Range("C6:E6").NumberFormat = "#,##0"

# --- Append two new company rows -------------------------------------------
# Shared strings must be introduced in a specific order so the shared
# string table indices line up with the target workbook (UENs first, then
# the two new company names, KUDOS before NUS).
$ws.Range("A7").Value = "200604346E"
$ws.Range("H7").Value = "200604346E"

$ws.Range("A8").Value = "201626142G"
$ws.Range("H8").Value = "201626142G"

$ws.Range("I8").Value = "KUDOS DATA PTE. LTD."
$ws.Range("I7").Value = "NATIONAL UNIVERSITY OF SINGAPORE"

$ws.Range("B7").Value = "S9000888C"
$ws.Range("C7").Value = 100000
$ws.Range("D7").Value = 100000000
$ws.Range("E7").Value = 14000
$ws.Range("F7").Value = 2008
$ws.Range("G7").Value = "Singapore"

$ws.Range("B8").Value = "S9000018G"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 10000
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 2016
$ws.Range("G8").Value = "Singapore"

$ws.Range("C7:E8").NumberFormat = "#,##0"
$ws.Range("A7:I8").Borders.LineStyle = 1
